# Update the MeasureReportAssignedPractitioner StructureDefinition workbook
# to the new published version (5.0.0 -> 6.0.0), per the
# "Alvearie/alvearie-fhir-ig" gh-pages deploy.

$wb = $excel.ActiveWorkbook

# ---- "Metadata" sheet ----
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# New publication date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$ws.Range("B9").Value = "Alvearie Team"

# The old "Contact" / "No display for ContactDetail" row is replaced with
# a "Jurisdiction" / "United States of America" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The second, duplicated "Contact" / "No display for ContactDetail" row is
# removed entirely - everything below shifts up by one row.
$ws.Rows.Item(11).Delete()

# ---- "Elements" sheet ----
$ws2 = $wb.Worksheets.Item("Elements")

# Root Extension row's Short/Definition now mirror the StructureDefinition's
# own Title/Description instead of the generic Extension placeholders.
$ws2.Range("K2").Value = "Measure Report Assigned Practitioner"
$ws2.Range("L2").Value = "The provider identifier of the physician(s) assigned to the rule measure using one of the physician attribution methods"
